$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 53; this shifts the existing rows 53-80
# down to 54-81 (which matches the diff's "renumbering" of every
# subsequent data row) and grows the used range to A1:R81.
$ws.Rows(53).Insert()

# Populate the newly inserted row 53 with the new weekly data point.
$ws.Cells.Item(53, 1).Value = 8
$ws.Cells.Item(53, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(53, 3).Value = "Coquimbo"
$ws.Cells.Item(53, 4).Value = 44455
$ws.Cells.Item(53, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(53, 5).Value = 4
$ws.Cells.Item(53, 6).Value = 100112001
$ws.Cells.Item(53, 7).Value = "Berenjena"
$ws.Cells.Item(53, 8).Value = "Sin especificar"
$ws.Cells.Item(53, 9).Value = "Primera"
$ws.Cells.Item(53, 10).Value = 580
$ws.Cells.Item(53, 11).Value = 9000
$ws.Cells.Item(53, 12).Value = 10000
$ws.Cells.Item(53, 13).Value = 9500
$ws.Cells.Item(53, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(53, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(53, 16).Value = 158
$ws.Cells.Item(53, 17).Value = 60
$ws.Cells.Item(53, 18).Value = "Hortaliza"
